$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 2 (pushes the existing "Nov 10 2020" row down to row 3)
$ws.Rows.Item(2).Insert()

# ---- New row 2: Oct 24 2020 vs Kolkata Knight Riders ----
$ws.Range("A2").Value = " Oct 24 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "KKR won by 59 runs"
$ws.Range("D2").Value = "Delhi Capitals"
$ws.Range("E2").Value = "Kolkata Knight Riders"
$ws.Range("F2").Value = "Ajinkya Rahane "
$ws.Range("G2").Formula = '="0"'
$ws.Range("H2").Formula = '="1"'
$ws.Range("I2").Formula = '="0"'
$ws.Range("J2").Formula = '="0"'
$ws.Range("K2").Formula = '="0.00"'

# ---- Row 4: Nov 2 2020 vs Royal Challengers Bangalore ----
$ws.Range("A4").Value = " Nov 2 2020"
$ws.Range("B4").Value = " Abu Dhabi"
$ws.Range("C4").Value = "Capitals won by 6 wickets (with 6 balls remaining)"
$ws.Range("D4").Value = "Delhi Capitals"
$ws.Range("E4").Value = "Royal Challengers Bangalore"
$ws.Range("F4").Value = "Ajinkya Rahane "
$ws.Range("G4").Formula = '="60"'
$ws.Range("H4").Formula = '="46"'
$ws.Range("I4").Formula = '="5"'
$ws.Range("J4").Formula = '="1"'
$ws.Range("K4").Formula = '="130.43"'

# ---- Row 5: Oct 11 2020 vs Mumbai Indians ----
$ws.Range("A5").Value = " Oct 11 2020"
$ws.Range("B5").Value = " Abu Dhabi"
$ws.Range("C5").Value = "Mumbai won by 5 wickets (with 2 balls remaining)"
$ws.Range("D5").Value = "Delhi Capitals"
$ws.Range("E5").Value = "Mumbai Indians"
$ws.Range("F5").Value = "Ajinkya Rahane "
$ws.Range("G5").Formula = '="15"'
$ws.Range("H5").Formula = '="15"'
$ws.Range("I5").Formula = '="3"'
$ws.Range("J5").Formula = '="0"'
$ws.Range("K5").Formula = '="100.00"'

# ---- Row 6: Nov 5 2020 vs Mumbai Indians ----
$ws.Range("A6").Value = " Nov 5 2020"
$ws.Range("B6").Value = " Dubai (DSC)"
$ws.Range("C6").Value = "Mumbai won by 57 runs"
$ws.Range("D6").Value = "Delhi Capitals"
$ws.Range("E6").Value = "Mumbai Indians"
$ws.Range("F6").Value = "Ajinkya Rahane "
$ws.Range("G6").Formula = '="0"'
$ws.Range("H6").Formula = '="3"'
$ws.Range("I6").Formula = '="0"'
$ws.Range("J6").Formula = '="0"'
$ws.Range("K6").Formula = '="0.00"'

# ---- Row 7: Oct 17 2020 vs Chennai Super Kings ----
$ws.Range("A7").Value = " Oct 17 2020"
$ws.Range("B7").Value = " Sharjah"
$ws.Range("C7").Value = "Capitals won by 5 wickets (with 1 ball remaining)"
$ws.Range("D7").Value = "Delhi Capitals"
$ws.Range("E7").Value = "Chennai Super Kings"
$ws.Range("F7").Value = "Ajinkya Rahane "
$ws.Range("G7").Formula = '="8"'
$ws.Range("H7").Formula = '="10"'
$ws.Range("I7").Formula = '="1"'
$ws.Range("J7").Formula = '="0"'
$ws.Range("K7").Formula = '="80.00"'

# ---- Row 8: Oct 14 2020 vs Rajasthan Royals ----
$ws.Range("A8").Value = " Oct 14 2020"
$ws.Range("B8").Value = " Dubai (DSC)"
$ws.Range("C8").Value = "Capitals won by 13 runs"
$ws.Range("D8").Value = "Delhi Capitals"
$ws.Range("E8").Value = "Rajasthan Royals"
$ws.Range("F8").Value = "Ajinkya Rahane "
$ws.Range("G8").Formula = '="2"'
$ws.Range("H8").Formula = '="9"'
$ws.Range("I8").Formula = '="0"'
$ws.Range("J8").Formula = '="0"'
$ws.Range("K8").Formula = '="22.22"'

# Convert the helper formulas (used to force text-typed numeric-looking values)
# into plain static text values, matching the original workbook's plain text cells.
# (Each row is copy/paste-special'd separately, since multi-area ranges only
# reliably paste-special their first area.)
foreach ($r in @(2, 4, 5, 6, 7, 8)) {
    $rowRange = $ws.Range("G" + $r + ":K" + $r)
    $rowRange.Copy()
    $rowRange.PasteSpecial(-4163)
}
$excel.CutCopyMode = 0
